# Update the "Förändrad" (Changed) date column (C) for all data rows
# from the old date serial 45192 to the new date serial 45202.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 218; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 45192) {
        $cell.Value2 = 45202
    }
}
